# Registration Data Import XLS Template - replace placeholder phone numbers
# on the "Individuals" sheet with realistic sample phone numbers.
#
# Column H = phone_number_1, Column I = phone_number_2.
# Rows 3-29 are the 27 sample individuals. The original placeholder data
# alternated row-by-row between two patterns; the new data keeps the same
# alternation but swaps in real-looking phone numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individuals")

for ($row = 3; $row -le 29; $row++) {
    if (($row % 2) -eq 1) {
        $ws.Cells.Item($row, 8).Value = "+44 1632 960852"
        $ws.Cells.Item($row, 9).Value = "+1-541-754-3010"
    } else {
        $ws.Cells.Item($row, 8).Value = "+1-613-555-0182"
        $ws.Cells.Item($row, 9).Value = "+36 55 979 922"
    }
}
